# Apply cell value updates to match the target crypto price snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text / safe (non-numeric-looking) values: direct assignment ---
$ws.Range("D2").Value = "27.287.43"
$ws.Range("E2").Value = "  -3.00%  "
$ws.Range("D3").Value = "1.853.72"
$ws.Range("E3").Value = "  -3.77%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("E5").Value = "  -1.69%  "
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("E7").Value = "  -3.92%  "
$ws.Range("E8").Value = "  -4.56%  "
$ws.Range("E9").Value = "  -8.86%  "
$ws.Range("E11").Value = "  -3.21%  "
$ws.Range("E12").Value = "  -4.27%  "
$ws.Range("D13").Value = "1.858.32"
$ws.Range("E13").Value = "  -3.93%  "
$ws.Range("E14").Value = "  -3.47%  "
$ws.Range("E15").Value = "  -5.21%  "
$ws.Range("E16").Value = "  +0.00%  "
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("E17").Value = "  +0.37%  "
$ws.Range("B18").Value = "Litecoin"
$ws.Range("C18").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("E18").Value = "  -5.25%  "
$ws.Range("E19").Value = "  -3.71%  "
$ws.Range("E20").Value = "  -5.55%  "
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("E22").Value = "  -4.52%  "
$ws.Range("D23").Value = "27.298.62"
$ws.Range("E23").Value = "  -3.00%  "
$ws.Range("E24").Value = "  -4.61%  "
$ws.Range("E25").Value = "  +0.27%  "
$ws.Range("D26").Value = "2.074.51"
$ws.Range("E26").Value = "  -4.08%  "
$ws.Range("E27").Value = "  -0.10%  "
$ws.Range("E28").Value = "  -1.06%  "
$ws.Range("E29").Value = "  -4.56%  "
$ws.Range("E30").Value = "  -4.42%  "
$ws.Range("E31").Value = "  -2.03%  "
$ws.Range("E32").Value = "  -2.88%  "
$ws.Range("E33").Value = "  -4.37%  "
$ws.Range("E34").Value = "  +0.16%  "
$ws.Range("E35").Value = "  -1.36%  "
$ws.Range("E36").Value = "  -5.55%  "
$ws.Range("E37").Value = "  -2.79%  "
$ws.Range("E38").Value = "  -4.11%  "
$ws.Range("E39").Value = "  -1.48%  "
$ws.Range("E40").Value = "  -11.26%  "
$ws.Range("E41").Value = "  -0.09%  "
$ws.Range("E42").Value = "  -4.28%  "
$ws.Range("E43").Value = "  -1.26%  "
$ws.Range("E44").Value = "  -8.30%  "
$ws.Range("E45").Value = "  -1.10%  "
$ws.Range("E46").Value = "  -4.85%  "
$ws.Range("E47").Value = "  -6.49%  "
$ws.Range("E48").Value = "  -2.80%  "
$ws.Range("E49").Value = "  -5.91%  "
$ws.Range("E50").Value = "  -1.31%  "
$ws.Range("E51").Value = "  -1.68%  "

# --- Numeric-looking text values in column D: must stay text (not be coerced to a number). ---
# Use a scratch cell formatted as Text, then Copy/PasteSpecial(values) into the target so the
# target cell keeps its original (default) style while the pasted value remains a text string.
$helper = $ws.Range("ZZ1")
$helper.NumberFormat = "@"
$helper.Value = "1.002"
$helper.Copy()
$ws.Range("D4").PasteSpecial(-4163)
$helper.Value = "323.90"
$helper.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$helper.Value = "0.4542"
$helper.Copy()
$ws.Range("D7").PasteSpecial(-4163)
$helper.Value = "0.3876"
$helper.Copy()
$ws.Range("D8").PasteSpecial(-4163)
$helper.Value = "48.27"
$helper.Copy()
$ws.Range("D9").PasteSpecial(-4163)
$helper.Value = "0.07921"
$helper.Copy()
$ws.Range("D10").PasteSpecial(-4163)
$helper.Value = "1.014"
$helper.Copy()
$ws.Range("D11").PasteSpecial(-4163)
$helper.Value = "21.35"
$helper.Copy()
$ws.Range("D12").PasteSpecial(-4163)
$helper.Value = "5.890"
$helper.Copy()
$ws.Range("D14").PasteSpecial(-4163)
$helper.Value = "7.125"
$helper.Copy()
$ws.Range("D15").PasteSpecial(-4163)
$helper.Value = "1.002"
$helper.Copy()
$ws.Range("D16").PasteSpecial(-4163)
$helper.Value = "0.06600"
$helper.Copy()
$ws.Range("D17").PasteSpecial(-4163)
$helper.Value = "85.77"
$helper.Copy()
$ws.Range("D18").PasteSpecial(-4163)
$helper.Value = "0.00001027"
$helper.Copy()
$ws.Range("D19").PasteSpecial(-4163)
$helper.Value = "17.08"
$helper.Copy()
$ws.Range("D20").PasteSpecial(-4163)
$helper.Value = "5.497"
$helper.Copy()
$ws.Range("D22").PasteSpecial(-4163)
$helper.Value = "10.88"
$helper.Copy()
$ws.Range("D24").PasteSpecial(-4163)
$helper.Value = "2.289"
$helper.Copy()
$ws.Range("D25").PasteSpecial(-4163)
$helper.Value = "153.87"
$helper.Copy()
$ws.Range("D27").PasteSpecial(-4163)
$helper.Value = "2.057"
$helper.Copy()
$ws.Range("D29").PasteSpecial(-4163)
$helper.Value = "5.462"
$helper.Copy()
$ws.Range("D30").PasteSpecial(-4163)
$helper.Value = "121.18"
$helper.Copy()
$ws.Range("D31").PasteSpecial(-4163)
$helper.Value = "0.09327"
$helper.Copy()
$ws.Range("D32").PasteSpecial(-4163)
$helper.Value = "0.9360"
$helper.Copy()
$ws.Range("D33").PasteSpecial(-4163)
$helper.Value = "1.453"
$helper.Copy()
$ws.Range("D34").PasteSpecial(-4163)
$helper.Value = "3.588"
$helper.Copy()
$ws.Range("D35").PasteSpecial(-4163)
$helper.Value = "5.254"
$helper.Copy()
$ws.Range("D36").PasteSpecial(-4163)
$helper.Value = "0.06011"
$helper.Copy()
$ws.Range("D37").PasteSpecial(-4163)
$helper.Value = "0.02225"
$helper.Copy()
$ws.Range("D38").PasteSpecial(-4163)
$helper.Value = "1.220"
$helper.Copy()
$ws.Range("D39").PasteSpecial(-4163)
$helper.Value = "8.043"
$helper.Copy()
$ws.Range("D40").PasteSpecial(-4163)
$helper.Value = "0.5913"
$helper.Copy()
$ws.Range("D42").PasteSpecial(-4163)
$helper.Value = "0.1883"
$helper.Copy()
$ws.Range("D43").PasteSpecial(-4163)
$helper.Value = "10.15"
$helper.Copy()
$ws.Range("D44").PasteSpecial(-4163)
$helper.Value = "0.5603"
$helper.Copy()
$ws.Range("D46").PasteSpecial(-4163)
$helper.Value = "12.01"
$helper.Copy()
$ws.Range("D47").PasteSpecial(-4163)
$helper.Value = "1.916"
$helper.Copy()
$ws.Range("D49").PasteSpecial(-4163)
$helper.Value = "0.06733"
$helper.Copy()
$ws.Range("D50").PasteSpecial(-4163)
$helper.Value = "108.32"
$helper.Copy()
$ws.Range("D51").PasteSpecial(-4163)
$helper.Clear()
$excel.CutCopyMode = $false
